$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(222, 'GORC International Model WG', 'Data Infrastructures and Environments - Institutional'),
    @(223, 'Global Open Research Commons IG', 'Data Infrastructures and Environments - Institutional'),
    @(224, 'National PID Strategies Interest Group', 'Data Infrastructures and Environments - International'),
    @(225, 'RDA-OfR Mapping the digital research data infrastructure landscape WG', 'Data Infrastructures and Environments - International'),
    @(226, 'RDA / CODATA Data Systems, Tools, and Services for Crisis Situations WG', 'Data Infrastructures and Environments - International'),
    @(227, 'Working with PIDs in Tools IG', 'Data Infrastructures and Environments - International'),
    @(228, 'Metadata IG', 'Data Infrastructures and Environments - International'),
    @(229, 'Data Repository Attributes WG', 'Data Infrastructures and Environments - International'),
    @(230, 'RDA/CODATA Materials Data, Infrastructure & Interoperability IG', 'Data Infrastructures and Environments - Regional or Disciplinary'),
    @(231, 'GORC International Model WG', 'Data Infrastructures and Environments - Regional or Disciplinary'),
    @(232, 'DMP Common Standards WG', 'Data Infrastructures and Environments - Regional or Disciplinary'),
    @(233, 'Science communication for research data', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(234, 'RDA-OfR Mapping the digital research data infrastructure landscape WG', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(235, 'Ethics and Social Aspects of Data IG', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(236, 'Reproducibility IG', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(237, 'Evaluation of Research IG', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(238, 'Working with PIDs in Tools IG', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(239, 'National PID Strategies Interest Group', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(240, 'Domain Repositories IG', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(241, 'Research Data Architectures in Research Institutions IG', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(242, 'Linguistics Data IG', 'Data Lifecycles - Versioning, Provenance, Citation, and Reward'),
    @(243, 'Science communication for research data', 'Discipline Focused Data Issues'),
    @(244, 'RDA-OfR Mapping the digital research data infrastructure landscape WG', 'Discipline Focused Data Issues'),
    @(245, 'ESIP/RDA Earth, Space, and Environmental Sciences IG', 'Discipline Focused Data Issues'),
    @(246, 'National PID Strategies Interest Group', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(247, 'RDA / CODATA Data Systems, Tools, and Services for Crisis Situations WG', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(248, 'CoreTrustSeal Maintenance WG', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(249, 'Working with PIDs in Tools IG', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(250, 'Domain Repositories IG', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(251, 'Artificial Intelligence and Data Visitation (AIDV) WG', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(252, 'Metadata IG', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(253, 'Research data needs of the Photon and Neutron Science community IG', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(254, 'Vocabulary Services IG', 'FAIR, CARE, TRUST - Adoption, Implementation, and Deployment'),
    @(255, 'RDA / CODATA Data Systems, Tools, and Services for Crisis Situations WG', 'FAIR, CARE, TRUST - Principles'),
    @(256, 'Ethics and Social Aspects of Data IG', 'FAIR, CARE, TRUST - Principles'),
    @(257, 'CoreTrustSeal Maintenance WG', 'FAIR, CARE, TRUST - Principles'),
    @(258, 'Artificial Intelligence and Data Visitation (AIDV) WG', 'FAIR, CARE, TRUST - Principles'),
    @(259, 'FAIR Digital Object Fabric IG', 'FAIR, CARE, TRUST - Principles'),
    @(260, 'Complex Citations Working Group', 'FAIR, CARE, TRUST - Principles, Semantics, Ontology, Standardisation'),
    @(261, 'Reproducibility IG', 'Other'),
    @(262, 'Evaluation of Research IG', 'Other'),
    @(263, 'CODATA/RDA Research Data Science Schools for Low and Middle Income Countries', 'Other'),
    @(264, 'Early Career and Engagement IG', 'Other'),
    @(265, 'RDA/WDS Scholarly Link Exchange (Scholix) WG', 'Research Software'),
    @(266, 'National PID Strategies Interest Group', 'Semantics, Ontology, Standardisation'),
    @(267, 'Research Data Management in Engineering IG', 'Semantics, Ontology, Standardisation'),
    @(268, 'RDA/CODATA Materials Data, Infrastructure & Interoperability IG', 'Semantics, Ontology, Standardisation'),
    @(269, 'Data Repository Attributes WG', 'Semantics, Ontology, Standardisation'),
    @(270, 'DMP Common Standards WG', 'Semantics, Ontology, Standardisation'),
    @(271, 'Ethics and Social Aspects of Data IG', 'Training, Stewardship, and Data Management Planning'),
    @(272, 'Linguistics Data IG', 'Training, Stewardship, and Data Management Planning,')
)

# Column B (Pathway) is populated first for all new rows, then column A (Group) -
# this matches the order new labels were introduced into the shared-strings table.
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
}

$ws.Range("A276").Select()

Write-Output "done"
